$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    # wdFindContinue=1, wdReplaceAll=2 -> replaces every matching occurrence
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Peserta 1 (kolom pertama)
Replace-Text "T1" "C1"
Replace-Text "SUNARTO" "ALFIAN SATYA"
Replace-Text "XXL" "L"

# Peserta 2 (kolom kedua)
Replace-Text "T2" "C2"
Replace-Text "NOOR MAULANA" "RAHMAN PALA"
Replace-Text "M" "XL"

# Peserta 3 (kolom ketiga)
Replace-Text "T3" "C3"
Replace-Text "OKTE DWI PANGGA" "DWI HENDRA A"
Replace-Text "XXL" "L"

# Field KELAS muncul 3x dengan teks lama/baru yang sama; ReplaceAll menangani semuanya
Replace-Text "DP4 NAUTIKA / 33" "DP 3 NAUTIKA/33"
